$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Comment on A1: the leading "Author:" run used to be bold Tahoma 9 (FF000000);
#        it now matches the plain Arial 10 font already used for this run on every other
#        header comment (B1..N1). ---
$cm = $ws.Range("A1").Comment
if ($cm -ne $null) {
    $tf = $cm.Shape.TextFrame
    $chars = $tf.Characters(1, 8)   # "Author:" + line feed
    $chars.Font.Name = "Arial"
    $chars.Font.Size = 10
    # Note: deliberately not touching Bold/FontStyle/Color here - on this host,
    # changing those on a comment's character run also re-stamps the parent
    # cell's own style, which would wrongly alter A1's (still-bold) cell font.
}

# --- 2. Column K ("Investment Date *") moves from a date number format to plain text,
#        right aligned - dates are now written as literal (local-parsed) strings instead
#        of Excel date serials. ---
$kCol = $ws.Range("K2:K3")
$kCol.NumberFormat = "@"
$kCol.HorizontalAlignment = -4152   # xlRight

# K2 held a real Excel date serial (displayed 12/05/25 under mm/dd/yy) -> literal text now.
$ws.Range("K2").Value = "12/05/2025"

# K3 was already literal text ("15/02/2025"); keep the value, just pick up the new style.
$ws.Range("K3").Value = "15/02/2025"

# --- 3. Saved selection moves to K3. ---
$ws.Range("K3").Select() | Out-Null
